$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore revision: cell C10 ("Integer min" for rule R30) changes from 18 to 1
$ws.Range("C10").Value = 1
